# Sample Project / Main.xlsx - "Rules" sheet
# Change cell B11 from the rule-name string "R40" to the text "1",
# keeping the cell's existing style/format (General, quote-prefix free).
#
# A plain `.Value = "1"` assignment on a General-formatted cell is
# interpreted by Excel as a number, and pre-setting NumberFormat="@" (or
# using a leading apostrophe) marks the cell with a quote-prefix and a new
# style. Neither matches the source change, which keeps the same style and
# stores a genuine text value. Routing the text through a formula result
# and then pasting it back as a value avoids both pitfalls.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Cells.Item(11, 2)   # B11
$cell.Formula = "=TEXT(1,""0"")"
$cell.Copy()
$cell.PasteSpecial(-4163)       # xlPasteValues
$excel.CutCopyMode = $false

$wb.Save()
